# Position sizing calculator: add a "Max Risk %" parameter cell (B23, labeled
# by B22 "MAX_RSIK") and rewire the per-row "Max Risk per Trade" column (H) to
# reference it instead of the hard-coded 0.1 / 0.02 constants. Also swap the
# week-2 ticker in row 5 (was CNXC) for XLF with its own price/size data, and
# relabel the G column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Position Sizing")

# --- Header row tweaks -------------------------------------------------------
$ws.Range("G1").Value = "Real risk (%)"

# --- New parameter block: B22 label, B23 value (10%) -----------------------
$ws.Range("B22").Value = "MAX_RSIK"
$ws.Range("B23").Value = 0.1

$ws.Range("H1").Formula = '="Max Risk per Trade (" & B23*100 & "%) ($)"'

# --- Rewire the "Max Risk per Trade" formulas to use $B$23 ------------------
$ws.Range("H2").Formula = "=B2*`$B`$23"
$ws.Range("H3").Formula = "=B3*`$B`$23"
$ws.Range("H4").Formula = "=B4*`$B`$23"
$ws.Range("H5").Formula = "=B5*`$B`$23"
$ws.Range("H6").Formula = "=B6*`$B`$23"
$ws.Range("H7").Formula = "=B7*`$B`$23"
$ws.Range("H8").Formula = "=B8*`$B`$23"
$ws.Range("H9").Formula = "=B9*`$B`$23"
$ws.Range("H10").Formula = "=B10*`$B`$23"

# --- Week 2 stock swap: CNXC -> XLF -----------------------------------------
$ws.Range("A5").Value = "XLF"
$ws.Range("B5").Value = 380
$ws.Range("C5").Value = 54.32

# --- Recalculate and restore the reported selection -------------------------
$excel.CalculateFull()
$ws.Range("D7").Select()
